# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.060.16"
$ws.Range("E2").Value = "  -1.13%  "

# Row 3
$ws.Range("D3").Value = "'2.549.55"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'579.23"
$ws.Range("E5").Value = "  +0.84%  "

# Row 6
$ws.Range("D6").Value = "'146.91"
$ws.Range("E6").Value = "  -1.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.97%  "

# Row 9
$ws.Range("E9").Value = "  -1.22%  "

# Row 10
$ws.Range("E10").Value = "  -4.41%  "

# Row 11
$ws.Range("E11").Value = "  -0.67%  "

# Row 12
$ws.Range("E12").Value = "  -0.95%  "

# Row 13
$ws.Range("D13").Value = "'27.20"
$ws.Range("E13").Value = "  -3.70%  "

# Row 14
$ws.Range("D14").Value = "'3.005.70"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15
$ws.Range("D15").Value = "'62.976.26"
$ws.Range("E15").Value = "  -1.08%  "

# Row 16
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").Value = "'2.546.33"
$ws.Range("E17").Value = "  -0.13%  "

# Row 18
$ws.Range("E18").Value = "  -1.94%  "

# Row 19
$ws.Range("D19").Value = "'335.70"
$ws.Range("E19").Value = "  -2.09%  "

# Row 20
$ws.Range("E20").Value = "  -0.93%  "

# Row 21
$ws.Range("E21").Value = "  -2.23%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").Value = "'65.31"
$ws.Range("E23").Value = "  -1.23%  "

# Row 24
$ws.Range("E24").Value = "  -0.24%  "

# Row 25
$ws.Range("D25").Value = "'1.61"
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("E27").Value = "  +3.84%  "

# Row 28
$ws.Range("E28").Value = "  -1.32%  "

# Row 29
$ws.Range("D29").Value = "'7.30"
$ws.Range("E29").Value = "  +2.89%  "

# Row 30
$ws.Range("B30").Value = "'PancakeSwap"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.89"
$ws.Range("E30").Value = "  +0.64%  "

# Row 31
$ws.Range("B31").Value = "'PEPE"
$ws.Range("C31").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0₃0812"
$ws.Range("E31").Value = "  -2.96%  "

# Row 32
$ws.Range("D32").Value = "'178.14"
$ws.Range("E32").Value = "  +0.30%  "

# Row 33
$ws.Range("E33").Value = "  -4.21%  "

# Row 34
$ws.Range("D34").Value = "'406.47"
$ws.Range("E34").Value = "  -4.19%  "

# Row 35
$ws.Range("B35").Value = "'PolygonEcosystemToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "'0.400"
$ws.Range("E35").Value = "  -1.61%  "

# Row 36
$ws.Range("B36").Value = "'EthereumClassic"
$ws.Range("C36").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'19.12"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("E37").Value = "  +0.03%  "

# Row 38
$ws.Range("D38").Value = "'4.33"
$ws.Range("E38").Value = "  -3.07%  "

# Row 39
$ws.Range("E39").Value = "  -1.36%  "

# Row 40
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("D41").Value = "'39.49"
$ws.Range("E41").Value = "  -2.44%  "

# Row 42
$ws.Range("D42").Value = "'151.13"
$ws.Range("E42").Value = "  -3.30%  "

# Row 43
$ws.Range("E43").Value = "  -1.35%  "

# Row 44
$ws.Range("E44").Value = "  -0.87%  "

# Row 45
$ws.Range("E45").Value = "  +0.59%  "

# Row 46
$ws.Range("D46").Value = "'0.602"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47
$ws.Range("E47").Value = "  -0.40%  "

# Row 48
$ws.Range("E48").Value = "  +2.35%  "

# Row 49
$ws.Range("D49").Value = "'18.28"
$ws.Range("E49").Value = "  -2.73%  "

# Row 50
$ws.Range("E50").Value = "  -8.45%  "

# Row 51
$ws.Range("D51").Value = "'11.30"
$ws.Range("E51").Value = "  +0.40%  "
